$p = $ppt.ActivePresentation

$oldText = "https://www.youtube.com/watch?v=qcY-uiEHhn0&list=PLv2a_5pNAko2Jl4Ks7V428ttvy-Fj4NKU&index=2"
$newText = "https://www.youtube.com/watch?v=_V3dqC80FHU&list=PLv2a_5pNAko2Jl4Ks7V428ttvy-Fj4NKU&index=3"

for ($i = 2; $i -le 10; $i++) {
    $s = $p.Slides.Item($i)
    foreach ($shape in $s.Shapes) {
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
            }
        }
    }
}
